$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# D1 calculation corrections to R script, block order correction to
# psychopy script.
#
# The "orStimulus" block (Insects / Flowers, row 3) and the
# "leftAttribute/rightAttribute" block (Negative / Positive, row 2) were
# swapped: the category-label block now belongs on row 2 (together with the
# long instruction text) and the attribute-label block moved down to row 3
# (together with the second instruction paragraph).
# ---------------------------------------------------------------------------

# Capture the original values before anything is overwritten.
$origB2 = $ws.Range("B2").Value2
$origC2 = $ws.Range("C2").Value2
$origD2 = $ws.Range("D2").Value2
$origE2 = $ws.Range("E2").Value2

$origB3 = $ws.Range("B3").Value2
$origC3 = $ws.Range("C3").Value2
$origD3 = $ws.Range("D3").Value2
$origE3 = $ws.Range("E3").Value2

# Carry the existing visual formatting (fill colour, number format,
# protection) along with the block it belongs to, before the source cells
# get re-formatted in the next step.
$ws.Range("D2:E2").Copy() | Out-Null
$ws.Range("D3:E3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B2:C2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Move the values to their corrected rows.
$ws.Range("B2").Value2 = $origB3
$ws.Range("C2").Value2 = $origC3
$ws.Range("D3").Value2 = $origD2
$ws.Range("E3").Value2 = $origE2

$ws.Range("B3").Value2 = $origB2
$ws.Range("C3").Value2 = $origC2
$ws.Range("D2").Value2 = $origD3
$ws.Range("E2").Value2 = $origE3

# The cells that used to hold the highlighted / fillable entries now just
# display blank placeholders, so drop the highlight fill and lock them back
# down (matching the plain, protected style used elsewhere on the sheet).
$ws.Range("F2").Copy() | Out-Null
$ws.Range("D2:E2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D2:E2").Locked = $true

$ws.Range("F2").Copy() | Out-Null
$ws.Range("B3:C3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B3:C3").NumberFormat = "@"
$ws.Range("B3:C3").Locked = $true

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet view: the frozen/scrolled top-left cell is reset back to the top of
# the sheet, and the active selection now spans the cells that were just
# rearranged (B3:C3, D2:E2, F2:F3).
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B3,C3,D2,E2,F2,F3").Select() | Out-Null
